# Auto-generated edit script applying the cryptos.xlsx value updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.227.40'
$ws.Range("E2").Value = '  -1.41%  '
$ws.Range("D3").Value = '2.300.32'
$ws.Range("E3").Value = '  -2.29%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '''317.59'
$ws.Range("E5").Value = '  +0.98%  '
$ws.Range("D6").Value = '''103.19'
$ws.Range("E6").Value = '  -5.12%  '
$ws.Range("D7").Value = '''0.631'
$ws.Range("E7").Value = '  -0.85%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '''0.610'
$ws.Range("E9").Value = '  -0.89%  '
$ws.Range("D10").Value = '''39.67'
$ws.Range("E10").Value = '  -2.51%  '
$ws.Range("D11").Value = '''0.0906'
$ws.Range("E11").Value = '  -2.46%  '
$ws.Range("D12").Value = '''8.34'
$ws.Range("E12").Value = '  -2.45%  '
$ws.Range("E13").Value = '  -0.39%  '
$ws.Range("D15").Value = '''15.33'
$ws.Range("E15").Value = '  -3.53%  '
$ws.Range("D16").Value = '2.648.28'
$ws.Range("E16").Value = '  -2.26%  '
$ws.Range("D17").Value = '2.293.79'
$ws.Range("E17").Value = '  -2.70%  '
$ws.Range("D18").Value = '42.331.53'
$ws.Range("E18").Value = '  -1.18%  '
$ws.Range("D19").Value = '''7.41'
$ws.Range("E19").Value = '  -2.73%  '
$ws.Range("E20").Value = '  -0.56%  '
$ws.Range("E21").Value = '  +2.00%  '
$ws.Range("D22").Value = '''73.62'
$ws.Range("E22").Value = '  -3.80%  '
$ws.Range("D23").Value = '''280.63'
$ws.Range("E23").Value = '  +3.73%  '
$ws.Range("D24").Value = '''11.39'
$ws.Range("E24").Value = '  +19.63%  '
$ws.Range("E25").Value = '  -2.92%  '
$ws.Range("E26").Value = '  -0.25%  '
$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").Value = '''3.98'
$ws.Range("E27").Value = '  +0.42%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").Value = '''10.86'
$ws.Range("E28").Value = '  -4.25%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '''2.40'
$ws.Range("E29").Value = '  +6.20%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = '''23.03'
$ws.Range("E30").Value = '  -2.25%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = '''36.14'
$ws.Range("E31").Value = '  -1.13%  '
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").Value = '''164.43'
$ws.Range("E32").Value = '  -1.56%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '''0.0874'
$ws.Range("E33").Value = '  -3.74%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '''5.89'
$ws.Range("E34").Value = '  -3.69%  '
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").Value = '''0.136'
$ws.Range("E35").Value = '  +3.56%  '
$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = '''2.64'
$ws.Range("E36").Value = '  -9.96%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '''0.114'
$ws.Range("E37").Value = '  -6.37%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '''4.59'
$ws.Range("E38").Value = '  -1.55%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.0349'
$ws.Range("E39").Value = '  -2.94%  '
$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").Value = '''3.76'
$ws.Range("E40").Value = '  -1.24%  '
$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").Value = '''2.77'
$ws.Range("E41").Value = '  +3.97%  '
$ws.Range("B42").Value = 'BitcoinSV'
$ws.Range("C42").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D42").Value = '''100.53'
$ws.Range("E42").Value = '  -5.02%  '
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = '''1.46'
$ws.Range("E43").Value = '  -2.40%  '
$ws.Range("B44").Value = 'MultiversX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D44").Value = '''69.16'
$ws.Range("E44").Value = '  -3.62%  '
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").Value = '''0.226'
$ws.Range("E45").Value = '  -4.79%  '
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").Value = '''1.00'
$ws.Range("E46").Value = '  -0.31%  '
$ws.Range("B47").Value = 'Celestia'
$ws.Range("C47").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D47").Value = '''12.00'
$ws.Range("E47").Value = '  -3.61%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '''111.92'
$ws.Range("E48").Value = '  -1.66%  '
$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D49").Value = '''77.06'
$ws.Range("E49").Value = '  -3.31%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").Value = '''8.93'
$ws.Range("E50").Value = '  -1.63%  '
$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").Value = '''5.29'
$ws.Range("E51").Value = '  -4.16%  '

# Reset number format/style on cells that needed a force-text apostrophe so no stray
# text-quote style sticks to them (matches original inlineStr cells with no explicit style).
$resetCells = @("D5","D6","D7","D9","D10","D11","D12","D15","D19","D22","D23","D24","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $resetCells) {
    $ws.Range($addr).Style = "Normal"
}
